$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update dimension will adjust automatically as Excel recalculates used range.

# Row 2: FAPs / Wnt5a / Ror2 / FAPs
$ws.Range("A2").Value = "FAPs"
$ws.Range("B2").Value = "Wnt5a"
$ws.Range("C2").Value = "Ror2"
$ws.Range("D2").Value = "FAPs"
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 7.514794999999999
$ws.Range("H2").Value = 22.544385
$ws.Range("I2").Value = 0.992147452492356
$ws.Range("J2").Value = 0.992147452492356
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 4.669265333333333
$ws.Range("N2").Value = 14.007796
$ws.Range("O2").Value = 0.8833113458668934
$ws.Range("P2").Value = 0.8833113458668933
$ws.Range("Q2").Value = 35.08857178060666
$ws.Range("R2").Value = 315.7971460254599
$ws.Range("S2").Value = 0.8763751015594327
$ws.Range("T2").Value = 0.8763751015594325

# Row 3: FAPs / Wnt5a / Ror2 / Neutro
$ws.Range("A3").Value = "FAPs"
$ws.Range("B3").Value = "Wnt5a"
$ws.Range("C3").Value = "Ror2"
$ws.Range("D3").Value = "Neutro"
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 7.514794999999999
$ws.Range("H3").Value = 22.544385
$ws.Range("I3").Value = 0.992147452492356
$ws.Range("J3").Value = 0.992147452492356
$ws.Range("K3").Value = 1
$ws.Range("L3").Value = 0.3333333333333333
$ws.Range("M3").Value = 0.017474
$ws.Range("N3").Value = 0.052422
$ws.Range("O3").Value = 0.003305655463074583
$ws.Range("P3").Value = 0.003305655463074583
$ws.Range("Q3").Value = 0.13131352783
$ws.Range("R3").Value = 1.18182175047
$ws.Range("S3").Value = 0.003279697646506887
$ws.Range("T3").Value = 0.003279697646506887

# Row 4: FAPs / Wnt5a / Ror2 / sCs
$ws.Range("A4").Value = "FAPs"
$ws.Range("B4").Value = "Wnt5a"
$ws.Range("C4").Value = "Ror2"
$ws.Range("D4").Value = "sCs"
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 7.514794999999999
$ws.Range("H4").Value = 22.544385
$ws.Range("I4").Value = 0.992147452492356
$ws.Range("J4").Value = 0.992147452492356
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 0.5993529999999999
$ws.Range("N4").Value = 1.798059
$ws.Range("O4").Value = 0.1133829986700321
$ws.Range("P4").Value = 0.1133829986700321
$ws.Range("Q4").Value = 4.504014927634999
$ws.Range("R4").Value = 40.536134348715
$ws.Range("S4").Value = 0.1124926532864165
$ws.Range("T4").Value = 0.1124926532864165

# Row 5: sCs / Wnt5a / Ror2 / FAPs
$ws.Range("A5").Value = "sCs"
$ws.Range("B5").Value = "Wnt5a"
$ws.Range("C5").Value = "Ror2"
$ws.Range("D5").Value = "FAPs"
$ws.Range("E5").Value = 1
$ws.Range("F5").Value = 0.3333333333333333
$ws.Range("G5").Value = 0.05947733333333333
$ws.Range("H5").Value = 0.178432
$ws.Range("I5").Value = 0.007852547507643968
$ws.Range("J5").Value = 0.00785254750764397
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 4.669265333333333
$ws.Range("N5").Value = 14.007796
$ws.Range("O5").Value = 0.8833113458668934
$ws.Range("P5").Value = 0.8833113458668933
$ws.Range("Q5").Value = 0.2777154506524444
$ws.Range("R5").Value = 2.499439055872
$ws.Range("S5").Value = 0.006936244307460713
$ws.Range("T5").Value = 0.006936244307460713

# Row 6: sCs / Wnt5a / Ror2 / Neutro
$ws.Range("A6").Value = "sCs"
$ws.Range("B6").Value = "Wnt5a"
$ws.Range("C6").Value = "Ror2"
$ws.Range("D6").Value = "Neutro"
$ws.Range("E6").Value = 1
$ws.Range("F6").Value = 0.3333333333333333
$ws.Range("G6").Value = 0.05947733333333333
$ws.Range("H6").Value = 0.178432
$ws.Range("I6").Value = 0.007852547507643968
$ws.Range("J6").Value = 0.00785254750764397
$ws.Range("K6").Value = 1
$ws.Range("L6").Value = 0.3333333333333333
$ws.Range("M6").Value = 0.017474
$ws.Range("N6").Value = 0.052422
$ws.Range("O6").Value = 0.003305655463074583
$ws.Range("P6").Value = 0.003305655463074583
$ws.Range("Q6").Value = 0.001039306922666667
$ws.Range("R6").Value = 0.009353762304
$ws.Range("S6").Value = [double]"2.595781656769598E-05"
$ws.Range("T6").Value = [double]"2.595781656769599E-05"

# Row 7: sCs / Wnt5a / Ror2 / sCs
$ws.Range("A7").Value = "sCs"
$ws.Range("B7").Value = "Wnt5a"
$ws.Range("C7").Value = "Ror2"
$ws.Range("D7").Value = "sCs"
$ws.Range("E7").Value = 1
$ws.Range("F7").Value = 0.3333333333333333
$ws.Range("G7").Value = 0.05947733333333333
$ws.Range("H7").Value = 0.178432
$ws.Range("I7").Value = 0.007852547507643968
$ws.Range("J7").Value = 0.00785254750764397
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 0.5993529999999999
$ws.Range("N7").Value = 1.798059
$ws.Range("O7").Value = 0.1133829986700321
$ws.Range("P7").Value = 0.1133829986700321
$ws.Range("Q7").Value = 0.03564791816533333
$ws.Range("R7").Value = 0.320831263488
$ws.Range("S7").Value = 0.0008903453836155596
$ws.Range("T7").Value = 0.0008903453836155597
